$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(2, "ECs", "Inhba", "Acvr2a", "ECs", 2, 0.6666666666666666, 4.123204333333334, 12.369613, 0.2909967288544799, 0.2909967288544799, 3, 1, 14.61878266666667, 43.856348, 0.2662829816142094, 0.2662829816142094, 60.27622803925823, 542.486052353324, 0.07748747659935257, 0.07748747659935257),
  @(3, "ECs", "Inhba", "Acvr2a", "FAPs", 2, 0.6666666666666666, 4.123204333333334, 12.369613, 0.2909967288544799, 0.2909967288544799, 3, 1, 27.084169, 81.25250700000001, 0.4933415757187404, 0.4933415757187404, 111.6735629855324, 1005.062066869791, 0.1435607847420682, 0.1435607847420682),
  @(4, "ECs", "Inhba", "Acvr2a", "sCs", 2, 0.6666666666666666, 4.123204333333334, 12.369613, 0.2909967288544799, 0.2909967288544799, 3, 1, 13.19647366666667, 39.589421, 0.2403754426670501, 0.2403754426670501, 54.41175740711923, 489.705816664073, 0.06994846751305916, 0.06994846751305916),
  @(5, "FAPs", "Inhba", "Acvr2a", "ECs", 3, 1, 8.433639666666666, 25.300919, 0.5952073574179045, 0.5952073574179045, 3, 1, 14.61878266666667, 43.856348, 0.2662829816142094, 0.2662829816142094, 123.2895453759791, 1109.605908383812, 0.158493589811954, 0.158493589811954),
  @(6, "FAPs", "Inhba", "Acvr2a", "FAPs", 3, 1, 8.433639666666666, 25.300919, 0.5952073574179045, 0.5952073574179045, 3, 1, 27.084169, 81.25250700000001, 0.4933415757187404, 0.4933415757187404, 228.4181220171037, 2055.763098153933, 0.2936405355879365, 0.2936405355879365),
  @(7, "FAPs", "Inhba", "Acvr2a", "sCs", 3, 1, 8.433639666666666, 25.300919, 0.5952073574179045, 0.5952073574179045, 3, 1, 13.19647366666667, 39.589421, 0.2403754426670501, 0.2403754426670501, 111.2943037753221, 1001.648733977899, 0.1430732320180139, 0.1430732320180139),
  @(8, "sCs", "Inhba", "Acvr2a", "ECs", 3, 1, 1.612402333333333, 4.837207, 0.1137959137276156, 0.1137959137276156, 3, 1, 14.61878266666667, 43.856348, 0.2662829816142094, 0.2662829816142094, 23.57135928222623, 212.142233540036, 0.03030191520290283, 0.03030191520290282),
  @(9, "sCs", "Inhba", "Acvr2a", "FAPs", 3, 1, 1.612402333333333, 4.837207, 0.1137959137276156, 0.1137959137276156, 3, 1, 27.084169, 81.25250700000001, 0.4933415757187404, 0.4933415757187404, 43.67057729199434, 393.0351956279491, 0.05614025538873572, 0.05614025538873571),
  @(10, "sCs", "Inhba", "Acvr2a", "sCs", 3, 1, 1.612402333333333, 4.837207, 0.1137959137276156, 0.1137959137276156, 3, 1, 13.19647366666667, 39.589421, 0.2403754426670501, 0.2403754426670501, 21.27802493190523, 191.502224387147, 0.02735374313597705, 0.02735374313597704)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
    $ws.Cells.Item($r, 15).Value = $row[15]
    $ws.Cells.Item($r, 16).Value = $row[16]
    $ws.Cells.Item($r, 17).Value = $row[17]
    $ws.Cells.Item($r, 18).Value = $row[18]
    $ws.Cells.Item($r, 19).Value = $row[19]
    $ws.Cells.Item($r, 20).Value = $row[20]
}
